$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# --- Swap columns B and C (values + number/cell formats travel together) ---
# Column B currently holds short codes (MAYRPEGR, LUISSADI, PAOLROSA) with a bordered style,
# column C holds the full names with the default style. The new layout puts the names in B
# (default style) and the codes in C (bordered style), so copy/paste the columns to swap them.
$ws.Range("B1:B3").Copy($ws.Range("E1:E3"))
$ws.Range("C1:C3").Copy($ws.Range("B1:B3"))
$ws.Range("E1:E3").Copy($ws.Range("C1:C3"))
$ws.Range("E1:E3").Clear()

# --- Update the employee name that changed ---
$ws.Range("B3").Value = "CAROLINA MARTINEZ MOLINA"

# --- Column A: repeat the (lower-cased) label on every row instead of a single merged cell ---
$ws.Range("A1:A3").UnMerge()
$ws.Range("A1").Value = "valor asegurado"
$ws.Range("A2").Value = "valor asegurado"
$ws.Range("A3").Value = "valor asegurado"

# Replace the old center-aligned / borderless look with a thin right border and no explicit
# horizontal alignment.
$ws.Range("A1:A3").HorizontalAlignment = 1
$ws.Range("A1:A3").Borders.Item(10).LineStyle = 1

# --- Column widths: widen column B to fit the names, drop column C's old custom width ---
$ws.Columns.Item(2).ColumnWidth = 41.3

# --- Selection as left by the editor ---
$ws.Range("C12").Select()
